$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the diff (ref -> new text value).
# Percentage-only values (e.g. "55%") are written via a NumberFormat
# toggle so Excel keeps them as literal text instead of converting
# them into a numeric percentage value.

$ws.Range("E2").Value = '2026-02-20 19:18:17'
$ws.Range("G2").Value = '212 cm'
$ws.Range("K2").Value = '11.7 MJ/m2'
$ws.Range("O2").Value = '0.7 °C'
$ws.Range("E3").Value = '2026-02-20 19:18:19'
$ws.Range("E4").Value = '2026-02-20 19:18:22'
$c = $ws.Range("H4")
$c.NumberFormat = "@"
$c.Value = '55%'
$c.NumberFormat = "General"
$ws.Range("J4").Value = '1022.3 hPa'
$ws.Range("O4").Value = '10.5 °C'
$ws.Range("E5").Value = '2026-02-20 19:18:24'
$ws.Range("N5").Value = '-6.0 °C 18:34 TU'
$ws.Range("E6").Value = '2026-02-20 19:18:26'
$c = $ws.Range("H6")
$c.NumberFormat = "@"
$c.Value = '67%'
$c.NumberFormat = "General"
$ws.Range("J6").Value = '1022.3 hPa'
$ws.Range("E7").Value = '2026-02-20 19:18:29'
$ws.Range("J7").Value = '1022.2 hPa'
$ws.Range("E8").Value = '2026-02-20 19:18:31'
$c = $ws.Range("H8")
$c.NumberFormat = "@"
$c.Value = '59%'
$c.NumberFormat = "General"
$ws.Range("O8").Value = '9.3 °C'
$ws.Range("E9").Value = '2026-02-20 19:18:34'
$ws.Range("E10").Value = '2026-02-20 19:18:36'
$c = $ws.Range("H10")
$c.NumberFormat = "@"
$c.Value = '78%'
$c.NumberFormat = "General"
$ws.Range("E11").Value = '2026-02-20 19:18:38'
$c = $ws.Range("H11")
$c.NumberFormat = "@"
$c.Value = '31%'
$c.NumberFormat = "General"
$ws.Range("E12").Value = '2026-02-20 19:18:41'
$ws.Range("E13").Value = '2026-02-20 19:18:43'
$ws.Range("J13").Value = '1023.2 hPa'
$ws.Range("E14").Value = '2026-02-20 19:18:45'
$c = $ws.Range("H14")
$c.NumberFormat = "@"
$c.Value = '55%'
$c.NumberFormat = "General"
$ws.Range("N14").Value = '8.5 °C 18:49 TU'
$ws.Range("O14").Value = '12.3 °C'
$ws.Range("E15").Value = '2026-02-20 19:18:48'
$ws.Range("E16").Value = '2026-02-20 19:18:50'
$ws.Range("M16").Value = '0.7 °C 18:39 TU'
$ws.Range("O16").Value = '-3.4 °C'
$ws.Range("E17").Value = '2026-02-20 19:18:52'
$ws.Range("O17").Value = '2.9 °C'
$ws.Range("E18").Value = '2026-02-20 19:18:55'
$ws.Range("J18").Value = '1022.6 hPa'
$ws.Range("E19").Value = '2026-02-20 19:18:57'
$ws.Range("O19").Value = '4.4 °C'
$ws.Range("E20").Value = '2026-02-20 19:18:59'
$ws.Range("E21").Value = '2026-02-20 19:19:01'
$ws.Range("J21").Value = '1022.2 hPa'
$ws.Range("O21").Value = '9.5 °C'
$ws.Range("E22").Value = '2026-02-20 19:19:04'
$ws.Range("O22").Value = '-4.2 °C'
$ws.Range("E23").Value = '2026-02-20 19:19:06'
$ws.Range("O23").Value = '-5.1 °C'
$ws.Range("E24").Value = '2026-02-20 19:19:09'
$ws.Range("O24").Value = '9.7 °C'
$ws.Range("E25").Value = '2026-02-20 19:19:11'
$ws.Range("O25").Value = '-1.6 °C'
$ws.Range("E26").Value = '2026-02-20 19:19:13'
$ws.Range("J26").Value = '1021.4 hPa'
$ws.Range("E27").Value = '2026-02-20 19:19:16'
$ws.Range("E28").Value = '2026-02-20 19:19:18'
$ws.Range("J28").Value = '1022.6 hPa'
$ws.Range("E29").Value = '2026-02-20 19:19:20'
$c = $ws.Range("H29")
$c.NumberFormat = "@"
$c.Value = '71%'
$c.NumberFormat = "General"
$ws.Range("O29").Value = '9.6 °C'
$ws.Range("E30").Value = '2026-02-20 19:19:23'
$ws.Range("J30").Value = '1022.0 hPa'
$ws.Range("E31").Value = '2026-02-20 19:19:25'
$ws.Range("J31").Value = '1021.2 hPa'
$ws.Range("E32").Value = '2026-02-20 19:19:27'
$ws.Range("E33").Value = '2026-02-20 19:19:30'
$ws.Range("J33").Value = '1022.6 hPa'
$ws.Range("O33").Value = '6.2 °C'
$ws.Range("E34").Value = '2026-02-20 19:19:32'
$ws.Range("M34").Value = '4.9 °C 18:32 TU'
$ws.Range("O34").Value = '0.6 °C'
$ws.Range("E35").Value = '2026-02-20 19:19:35'
$ws.Range("O35").Value = '4.0 °C'
$ws.Range("E36").Value = '2026-02-20 19:19:37'
$c = $ws.Range("H36")
$c.NumberFormat = "@"
$c.Value = '42%'
$c.NumberFormat = "General"
$ws.Range("J36").Value = '1022.2 hPa'
$ws.Range("E37").Value = '2026-02-20 19:19:39'
$c = $ws.Range("H37")
$c.NumberFormat = "@"
$c.Value = '65%'
$c.NumberFormat = "General"
$ws.Range("J37").Value = '1024.1 hPa'
$ws.Range("O37").Value = '5.0 °C'
$ws.Range("E38").Value = '2026-02-20 19:19:42'
$c = $ws.Range("H38")
$c.NumberFormat = "@"
$c.Value = '67%'
$c.NumberFormat = "General"
$ws.Range("E39").Value = '2026-02-20 19:19:44'
$ws.Range("O39").Value = '-2.8 °C'
$ws.Range("E40").Value = '2026-02-20 19:19:47'
$c = $ws.Range("H40")
$c.NumberFormat = "@"
$c.Value = '36%'
$c.NumberFormat = "General"
$ws.Range("J40").Value = '1023.1 hPa'
$ws.Range("E41").Value = '2026-02-20 19:19:49'
$ws.Range("J41").Value = '1022.8 hPa'
$ws.Range("E42").Value = '2026-02-20 19:19:51'
$c = $ws.Range("H42")
$c.NumberFormat = "@"
$c.Value = '66%'
$c.NumberFormat = "General"
$ws.Range("O42").Value = '10.2 °C'
$ws.Range("E43").Value = '2026-02-20 19:19:54'
$ws.Range("E44").Value = '2026-02-20 19:19:56'
$c = $ws.Range("H44")
$c.NumberFormat = "@"
$c.Value = '79%'
$c.NumberFormat = "General"
$ws.Range("M44").Value = '-1.1 °C 18:47 TU'
$ws.Range("O44").Value = '-4.9 °C'
$ws.Range("E45").Value = '2026-02-20 19:19:58'
$ws.Range("J45").Value = '1029.3 hPa'
$ws.Range("K45").Value = '8.7 MJ/m2'
$ws.Range("N45").Value = '1.5 °C 18:59 TU'
$ws.Range("O45").Value = '3.7 °C'
$ws.Range("E46").Value = '2026-02-20 19:20:01'
$ws.Range("J46").Value = '1026.2 hPa'
$ws.Range("K46").Value = '12.4 MJ/m2'
